$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New test case row 17 - "Customercare016"
$ws.Range("A17").Value = "Customercare016"
$ws.Range("B17").Value = "OPQA-5191||OPQA-5194||OPQA-5196||OPQA-5197"
$ws.Range("C17").Value = "Verify that the text `"Drug Research Advisor Customer Care `" should be hyperlinked and it should be linked to customer care / support page.||Verify that 'DRA_support@thomsonreuters.com' is replaced with ' Drug Research Advisor Customer Care.'||Verify that hyperlinked text `"Drug Research Advisor Customer Care `" should be opened in a second window / tab (based on user/browser preference)||Verify that the customer care page URL content shall be specific to DRA(Target Druggability)"
$ws.Range("D17").Value = "Y"

# Match formatting of the row above (Customercare015) for columns A, B, D
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)

# Description column needs top-aligned wrapped text (vs. center) for both rows
$ws.Range("C16:C17").VerticalAlignment = -4160
$ws.Range("C16:C17").WrapText = $true
$ws.Range("B17").VerticalAlignment = -4160
$ws.Range("B17").WrapText = $true

$ws.Rows.Item(17).RowHeight = 60

$ws.Range("C17").Select()
